$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value (45204, i.e. 2023-10-05)
# in every data row (rows 2 through 294). Bump it to 45205 (2023-10-06) for
# all of them - this mirrors the automatic "last updated" timestamp refresh.
$lastRow = 294
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45205
}
